$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rename: pipe_length -> length_pipe
$ws.Range("H1").Value = "length_pipe"

# Row 2 (concentration_soil)
$ws.Range("B2").Value = 0.0000355513009972861
$ws.Range("G2").Value = 3.468721592776412
$ws.Range("K2").Value = 0.0196

# Row 3 (pipe_length -> length_pipe)
$ws.Range("A3").Value = "length_pipe"
$ws.Range("B3").Value = 0.00003519930706972955
$ws.Range("G3").Value = 3.468721592776412
$ws.Range("K3").Value = 0.0196

# Row 4 (length_fraction_middle_point)
$ws.Range("B4").Value = 0.00003519930706972955
$ws.Range("G4").Value = 3.468721592776412
$ws.Range("K4").Value = 0.0196

# Row 5 (length_plume)
$ws.Range("B5").Value = 0.00003519930706972955
$ws.Range("G5").Value = 3.503408808704176
$ws.Range("K5").Value = 0.0196

# Row 6 (inner_diameter)
$ws.Range("B6").Value = 0.00003485079907894014
$ws.Range("G6").Value = 3.468721592776412
$ws.Range("K6").Value = 0.019796

# Row 7 (flow_rate)
$ws.Range("B7").Value = 0.00003519930706972955
$ws.Range("G7").Value = 3.468721592776412
$ws.Range("K7").Value = 0.0196

# Row 8 (log_Dp_ref)
$ws.Range("B8").Value = 0.0000354213398352793
$ws.Range("G8").Value = 3.468721592776412
$ws.Range("K8").Value = 0.0196

# Row 9 (log_Kpw_ref)
$ws.Range("B9").Value = 0.00003529185460926674
$ws.Range("G9").Value = 3.468721592776412
$ws.Range("K9").Value = 0.0196

# Row 10 (DIFFUSION_A_C)
$ws.Range("B10").Value = 0.00003499908219029517
$ws.Range("G10").Value = 3.468721592776412
$ws.Range("K10").Value = 0.0196

# Row 11 (PARTITIONING_A_C)
$ws.Range("B11").Value = 0.00003517707796672573
$ws.Range("G11").Value = 3.468721592776412
$ws.Range("K11").Value = 0.0196

# Row 12 (activattion_energy)
$ws.Range("B12").Value = 0.0000350433929244036
$ws.Range("G12").Value = 3.468721592776412
$ws.Range("K12").Value = 0.0196

# Row 13 (partitioning_enthalpie)
$ws.Range("B13").Value = 0.00003518401660445891
$ws.Range("G13").Value = 3.468721592776412
$ws.Range("K13").Value = 0.0196
